$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a ..." timestamp banner in A1 moves from 14:22 to 14:52
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 14:52"

# Daily refresh of the provincias/ciudades COVID table.
# A few rows swap which province-label they carry (the underlying data for
# that rank position changed), and a handful of rows just get refreshed
# case counts. Each block below sets the label (column A) together with
# Casos totales / Casos activos / Recuperados / Muertes (columns B-E).

# Row 18: was "A Coruña" -> now "Malaga"
$ws.Range("A18").Value = "Malaga"
$ws.Range("B18").Value = 2006
$ws.Range("C18").Value = 318
$ws.Range("D18").Value = 1548
$ws.Range("E18").Value = 140

# Row 19: was "Malaga" -> now "A Coruña"
$ws.Range("A19").Value = "A Coruña"
$ws.Range("B19").Value = 1969
$ws.Range("C19").Value = 333
$ws.Range("D19").Value = 1788
$ws.Range("E19").Value = 67

# Row 21: "Sevilla" - refreshed counts only
$ws.Range("B21").Value = 1757
$ws.Range("C21").Value = 158
$ws.Range("D21").Value = 1468
$ws.Range("E21").Value = 131

# Row 25: was "Cantabria" -> now "Granada"
$ws.Range("A25").Value = "Granada"
$ws.Range("B25").Value = 1600
$ws.Range("C25").Value = 177
$ws.Range("D25").Value = 1285
$ws.Range("E25").Value = 138

# Row 26: was "Granada" -> now "Cantabria"
$ws.Range("A26").Value = "Cantabria"
$ws.Range("B26").Value = 1572
$ws.Range("C26").Value = 175
$ws.Range("D26").Value = 1305
$ws.Range("E26").Value = 92

# Row 33: "Cordoba" - refreshed counts only
$ws.Range("B33").Value = 1069
$ws.Range("C33").Value = 101
$ws.Range("D33").Value = 926
$ws.Range("E33").Value = 42

# Row 35: was "Guadalajara" -> now "Jaen"
$ws.Range("A35").Value = "Jaen"
$ws.Range("B35").Value = 990
$ws.Range("C35").Value = 74
$ws.Range("D35").Value = 844
$ws.Range("E35").Value = 72

# Row 36: was "Jaen" -> now "Guadalajara"
$ws.Range("A36").Value = "Guadalajara"
$ws.Range("B36").Value = 973
$ws.Range("C36").Value = 1557
$ws.Range("D36").Value = 8976
$ws.Range("E36").Value = 133

# Row 38: was "Castello/Castellon" -> now "Cadiz"
$ws.Range("A38").Value = "Cadiz"
$ws.Range("B38").Value = 901
$ws.Range("C38").Value = 118
$ws.Range("D38").Value = 747
$ws.Range("E38").Value = 36

# Row 39: was "Cadiz" -> now "Castello/Castellon"
$ws.Range("A39").Value = "Castello/Castellon"
$ws.Range("B39").Value = 899
$ws.Range("C39").Value = 142
$ws.Range("D39").Value = 668
$ws.Range("E39").Value = 89

# Row 50: "Almeria" - refreshed counts only
$ws.Range("B50").Value = 378
$ws.Range("C50").Value = 54
$ws.Range("D50").Value = 298
$ws.Range("E50").Value = 26

# Row 52: "Huelva" - refreshed counts only
$ws.Range("B52").Value = 296
$ws.Range("C52").Value = 41
$ws.Range("D52").Value = 235
$ws.Range("E52").Value = 20
